$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): I1 = "I0", J1 = "IF" with same style/formatting as other headers (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data values for columns I (9) and J (10), rows 2-24
$data = @(
    @{ Row = 2;  I = 5; J = 5 },
    @{ Row = 3;  I = 6; J = 8 },
    @{ Row = 4;  I = 7; J = 9 },
    @{ Row = 5;  I = 8; J = 8 },
    @{ Row = 6;  I = 6; J = 7 },
    @{ Row = 7;  I = 6; J = 7 },
    @{ Row = 8;  I = 5; J = 6 },
    @{ Row = 9;  I = 1; J = 1 },
    @{ Row = 10; I = 1; J = 4 },
    @{ Row = 11; I = 6; J = 7 },
    @{ Row = 12; I = 1; J = 5 },
    @{ Row = 13; I = 1; J = 7 },
    @{ Row = 14; I = 1; J = 5 },
    @{ Row = 15; I = 1; J = 6 },
    @{ Row = 16; I = 1; J = 6 },
    @{ Row = 17; I = 1; J = 7 },
    @{ Row = 18; I = 1; J = 6 },
    @{ Row = 19; I = 1; J = 4 },
    @{ Row = 20; I = 1; J = 6 },
    @{ Row = 21; I = 1; J = 4 },
    @{ Row = 22; I = 6; J = 8 },
    @{ Row = 23; I = 1; J = 2 },
    @{ Row = 24; I = 3; J = 3 }
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 9).Value = $entry.I
    $ws.Cells.Item($entry.Row, 10).Value = $entry.J
}
